$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.440.27'
$ws.Range('E2').Value = '  -2.46%  '
$ws.Range('D3').Value = '2.889.89'
$ws.Range('E3').Value = '  -2.06%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range('E4').Value = '  +0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '567.76'
$c.ClearFormats()
$ws.Range('E5').Value = '  -4.56%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '143.97'
$c.ClearFormats()
$ws.Range('E6').Value = '  -2.98%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('D9').Value = '2.888.83'
$ws.Range('E9').Value = '  -2.06%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '7.00'
$c.ClearFormats()
$ws.Range('E10').Value = '  -3.82%  '
$ws.Range('E11').Value = '  -2.69%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.431'
$c.ClearFormats()
$ws.Range('E12').Value = '  -2.44%  '
$ws.Range('E13').Value = '  -1.59%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '32.04'
$c.ClearFormats()
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').Value = '3.368.89'
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('D17').Value = '61.457.80'
$ws.Range('E17').Value = '  -2.29%  '
$ws.Range('E18').Value = '  -2.50%  '
$ws.Range('D19').Value = '2.886.25'
$ws.Range('E19').Value = '  -2.45%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '432.08'
$c.ClearFormats()
$ws.Range('E20').Value = '  -2.41%  '
$ws.Range('E21').Value = '  -2.54%  '
$ws.Range('E22').Value = '  -1.74%  '
$ws.Range('E23').Value = '  -3.03%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '79.29'
$c.ClearFormats()
$ws.Range('E24').Value = '  -2.19%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '12.01'
$c.ClearFormats()
$ws.Range('E25').Value = '  +1.64%  '
$ws.Range('E26').Value = '  +0.02%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '10.00'
$c.ClearFormats()
$ws.Range('E27').Value = '  -11.43%  '
$ws.Range('E28').Value = '  -6.13%  '
$ws.Range('E29').Value = '  +1.68%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '6.99'
$c.ClearFormats()
$ws.Range('E30').Value = '  -2.98%  '
$ws.Range('E31').Value = '  -4.41%  '
$ws.Range('E32').Value = '  -6.82%  '
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('E34').Value = '  -2.25%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '25.49'
$c.ClearFormats()
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.960'
$c.ClearFormats()
$ws.Range('E36').Value = '  -3.25%  '
$ws.Range('E37').Value = '  -3.80%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '48.88'
$c.ClearFormats()
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('E39').Value = '  -5.40%  '
$ws.Range('E40').Value = '  -10.50%  '
$ws.Range('E41').Value = '  -3.52%  '
$ws.Range('E42').Value = '  -2.76%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '39.03'
$c.ClearFormats()
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('E44').Value = '  -5.15%  '
$ws.Range('D45').Value = '2.709.46'
$ws.Range('E45').Value = '  +0.55%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '133.02'
$c.ClearFormats()
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('E47').Value = '  -0.56%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '344.05'
$c.ClearFormats()
$ws.Range('E48').Value = '  -5.14%  '
$ws.Range('E50').Value = '  -1.37%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '21.57'
$c.ClearFormats()
$ws.Range('E51').Value = '  -5.29%  '
